$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Marco Sartorelli "
$ws.Range("B7").Value = "Elia Barozzi | I Magnifici"
$ws.Range("C7").Value = "Riccardo Zeni | demobusters"
$ws.Range("D7").Value = "Edoardo Pomarolli | Modium"
$ws.Range("E7").Value = "Michele Merighi | Clitoriders"
$ws.Range("F7").Value = "Moris Benedetti | Gli Introvabili"
